$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.226.58"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.555.94"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'234.94"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "'649.06"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("D8").Value = "'0.396"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'0.988"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").Value = "3.552.22"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'0.201"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "'41.97"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "'6.49"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "4.226.88"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "94.698.47"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "'0.0000252"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "3.560.84"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'7.87"
$ws.Range("E19").Value = "  -5.67%  "
$ws.Range("D20").Value = "'12.60"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'17.67"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "'3.44"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'502.63"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("E24").Value = "  -7.54%  "
$ws.Range("D25").Value = "'0.0000192"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "'6.54"
$ws.Range("E26").Value = "  -3.72%  "
$ws.Range("D27").Value = "'94.53"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").Value = "3.747.89"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'12.36"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'3.01"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'11.28"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'0.139"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'0.177"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "'31.72"
$ws.Range("E36").Value = "  +4.86%  "
$ws.Range("D37").Value = "'0.554"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").Value = "'8.22"
$ws.Range("E38").Value = "  +6.64%  "
$ws.Range("D39").Value = "'561.46"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").Value = "'1.53"
$ws.Range("E40").Value = "  +5.17%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.149"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").Value = "'0.895"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("D44").Value = "'1.75"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'2.31"
$ws.Range("E45").Value = "  +5.66%  "
$ws.Range("D46").Value = "'33.87"
$ws.Range("E46").Value = "  +31.39%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'23.46"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "'5.59"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "'3.60"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "'0.0408"
$ws.Range("E50").Value = "  -4.61%  "
$ws.Range("D51").Value = "'53.11"
$ws.Range("E51").Value = "  -1.56%  "
